$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.008688450532821
$ws.Range("D2").Value = 1.011347374404128
$ws.Range("E2").Value = 1.011229710712755
$ws.Range("F2").Value = 1.01232709078867
$ws.Range("I2").Value = 1.022830103444465
$ws.Range("J2").Value = 1.013952826990632
$ws.Range("K2").Value = 1.014214819854
$ws.Range("L2").Value = 1.014097511185495
$ws.Range("M2").Value = 1.015191583664154
$ws.Range("N2").Value = 1.008830506242373

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.009948203634625
$ws.Range("D3").Value = 1.012449680323405
$ws.Range("E3").Value = 1.012302981766324
$ws.Range("F3").Value = 1.014899176263818
$ws.Range("I3").Value = 1.023104481063391
$ws.Range("J3").Value = 1.014842956478751
$ws.Range("K3").Value = 1.015120565814835
$ws.Range("L3").Value = 1.014974277492985
$ws.Range("M3").Value = 1.017563230859071
$ws.Range("N3").Value = 1.009122870040444

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01075578689404
$ws.Range("D4").Value = 1.013156291742462
$ws.Range("E4").Value = 1.012991344094922
$ws.Range("F4").Value = 1.01651518665792
$ws.Range("I4").Value = 1.023268875536985
$ws.Range("J4").Value = 1.015410997224989
$ws.Range("K4").Value = 1.015699432483559
$ws.Range("L4").Value = 1.015534922614298
$ws.Range("M4").Value = 1.01904944613505
$ws.Range("N4").Value = 1.009309434557271

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.011093513686318
$ws.Range("D5").Value = 1.013451783703671
$ws.Range("E5").Value = 1.013279291445182
$ws.Range("F5").Value = 1.017183157210875
$ws.Range("I5").Value = 1.023334863415529
$ws.Range("J5").Value = 1.015647928555391
$ws.Range("K5").Value = 1.015941086237034
$ws.Range("L5").Value = 1.015769041472887
$ws.Range("M5").Value = 1.019662819575479
$ws.Range("N5").Value = 1.009387248959039

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.011150115856857
$ws.Range("D6").Value = 1.013501306897168
$ws.Range("E6").Value = 1.013327555286461
$ws.Range("F6").Value = 1.017294648422417
$ws.Range("I6").Value = 1.023345760548856
$ws.Range("J6").Value = 1.015687601258604
$ws.Range("K6").Value = 1.015981561856658
$ws.Range("L6").Value = 1.015808259265182
$ws.Range("M6").Value = 1.019765141637691
$ws.Range("N6").Value = 1.009400278377643

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.010760306577285
$ws.Range("D7").Value = 1.013160246246495
$ws.Range("E7").Value = 1.012995197289315
$ws.Range("F7").Value = 1.016524156698481
$ws.Range("I7").Value = 1.023269769513482
$ws.Range("J7").Value = 1.015414170438911
$ws.Range("K7").Value = 1.015702668127805
$ws.Range("L7").Value = 1.015538057085926
$ws.Range("M7").Value = 1.019057686779217
$ws.Range("N7").Value = 1.009310476731673

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.00911577363819
$ws.Range("D8").Value = 1.011721297590554
$ws.Range("E8").Value = 1.01159370895715
$ws.Range("F8").Value = 1.013206441962646
$ws.Range("I8").Value = 1.022925567630325
$ws.Range("J8").Value = 1.01425530918712
$ws.Range("K8").Value = 1.014522429533401
$ws.Range("L8").Value = 1.014395216334709
$ws.Range("M8").Value = 1.016003211375736
$ws.Range("N8").Value = 1.008929858640762

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.00615858568265
$ws.Range("D9").Value = 1.009133483155734
$ws.Range("E9").Value = 1.009076105162491
$ws.Range("F9").Value = 1.006982351546846
$ws.Range("I9").Value = 1.022217169405401
$ws.Range("J9").Value = 1.012151225669846
$ws.Range("K9").Value = 1.012386268034593
$ws.Range("L9").Value = 1.012329088347824
$ws.Range("M9").Value = 1.010242588057758
$ws.Range("N9").Value = 1.008238720883246

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.004145138349328
$ws.Range("D10").Value = 1.007371301044737
$ws.Range("E10").Value = 1.007363671239686
$ws.Range("F10").Value = 1.002566916081195
$ws.Range("I10").Value = 1.02167456551055
$ws.Range("J10").Value = 1.010704881919846
$ws.Range("K10").Value = 1.010922404985722
$ws.Range("L10").Value = 1.01091480418429
$ws.Range("M10").Value = 1.006136380127135
$ws.Range("N10").Value = 1.007763588448497

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.003262818469038
$ws.Range("D11").Value = 1.006599028032638
$ws.Range("E11").Value = 1.006613665948987
$ws.Range("F11").Value = 1.00058893835483
$ws.Range("I11").Value = 1.021422462053551
$ws.Range("J11").Value = 1.01006778061674
$ws.Range("K11").Value = 1.0102786611112
$ws.Range("L11").Value = 1.010293241184894
$ws.Range("M11").Value = 1.004292479803728
$ws.Range("N11").Value = 1.007554286559971

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.002933462900034
$ws.Range("D12").Value = 1.006310741291837
$ws.Range("E12").Value = 1.006333762215631
$ws.Range("F12").Value = 0.9998440314039331
$ws.Range("I12").Value = 1.021326198589508
$ws.Range("J12").Value = 1.009829462773604
$ws.Range("K12").Value = 1.010038020133945
$ws.Range("L12").Value = 1.010060948802073
$ws.Range("M12").Value = 1.003597412787125
$ws.Range("N12").Value = 1.007475992353871

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.003004185015604
$ws.Range("D13").Value = 1.006372645207232
$ws.Range("E13").Value = 1.00639386284455
$ws.Range("F13").Value = 1.000004282531593
$ws.Range("I13").Value = 1.021346966781681
$ws.Range("J13").Value = 1.009880659068163
$ws.Range("K13").Value = 1.010089708180945
$ws.Range("L13").Value = 1.010110841043771
$ws.Range("M13").Value = 1.003746971202148
$ws.Range("N13").Value = 1.007492811864144

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.00323562722319
$ws.Range("D14").Value = 1.006575227591048
$ws.Range("E14").Value = 1.006590556153115
$ws.Range("F14").Value = 1.000527573789866
$ws.Range("I14").Value = 1.021414558620614
$ws.Range("J14").Value = 1.010048115521325
$ws.Range("K14").Value = 1.010258801027929
$ws.Range("L14").Value = 1.010274068947409
$ws.Range("M14").Value = 1.004235234193015
$ws.Range("N14").Value = 1.007547826046079

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.00337800990943
$ws.Range("D15").Value = 1.006699854455835
$ws.Range("E15").Value = 1.006711569458924
$ws.Range("F15").Value = 1.000848631014647
$ws.Range("I15").Value = 1.02145585549641
$ws.Range("J15").Value = 1.010151068312811
$ws.Range("K15").Value = 1.010362781242423
$ws.Range("L15").Value = 1.010374450182244
$ws.Range("M15").Value = 1.004534714655489
$ws.Range("N15").Value = 1.007581648751371

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.0042034701378
$ws.Range("D16").Value = 1.007422356096577
$ws.Range("E16").Value = 1.007413263983092
$ws.Range("F16").Value = 1.002696773201763
$ws.Range("I16").Value = 1.021690931723832
$ws.Range("J16").Value = 1.010746932500902
$ws.Range("K16").Value = 1.010964916496052
$ws.Range("L16").Value = 1.01095585885445
$ws.Range("M16").Value = 1.006257343655341
$ws.Range("N16").Value = 1.007777402782195

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.004718420055279
$ws.Range("D17").Value = 1.007873060300089
$ws.Range("E17").Value = 1.007851112730732
$ws.Range("F17").Value = 1.003838189946289
$ws.Range("I17").Value = 1.021833766633006
$ws.Range("J17").Value = 1.011117774275572
$ws.Range("K17").Value = 1.01133994621156
$ws.Range("L17").Value = 1.011318079923213
$ws.Range("M17").Value = 1.007320084782005
$ws.Range("N17").Value = 1.007899229453754

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.005017773465819
$ws.Range("D18").Value = 1.008135060266073
$ws.Range("E18").Value = 1.008105684228284
$ws.Range("F18").Value = 1.004497604098853
$ws.Range("I18").Value = 1.021915428028849
$ws.Range("J18").Value = 1.011333038674078
$ws.Range("K18").Value = 1.0115577442124
$ws.Range("L18").Value = 1.011528475460757
$ws.Range("M18").Value = 1.00793362539117
$ws.Range("N18").Value = 1.00796994580482

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.005119675635105
$ws.Range("D19").Value = 1.008224246077959
$ws.Range("E19").Value = 1.00819234890344
$ws.Range("F19").Value = 1.004721377347333
$ws.Range("I19").Value = 1.021942993562885
$ws.Range("J19").Value = 1.01140626292194
$ws.Range("K19").Value = 1.011631847690218
$ws.Range("L19").Value = 1.011600066452402
$ws.Range("M19").Value = 1.008141759450388
$ws.Range("N19").Value = 1.007994000479776

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.004663275416191
$ws.Range("D20").Value = 1.007824796170408
$ws.Range("E20").Value = 1.007804220620989
$ws.Range("F20").Value = 1.003716385962054
$ws.Range("I20").Value = 1.02181861295899
$ws.Range("J20").Value = 1.01107809447978
$ws.Range("K20").Value = 1.011299807675259
$ws.Range("L20").Value = 1.011279308519776
$ws.Range("M20").Value = 1.007206720153434
$ws.Range("N20").Value = 1.007886194192515

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.003167518453291
$ws.Range("D21").Value = 1.006515611962968
$ws.Range("E21").Value = 1.006532671602976
$ws.Range("F21").Value = 1.000373761351914
$ws.Range("I21").Value = 1.021394727233023
$ws.Range("J21").Value = 1.009998850237339
$ws.Range("K21").Value = 1.010209049864958
$ws.Range("L21").Value = 1.010226041816593
$ws.Range("M21").Value = 1.004091735580907
$ws.Range("N21").Value = 1.007531641048096

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.002217662525463
$ws.Range("D22").Value = 1.005684178328557
$ws.Range("E22").Value = 1.005725548895555
$ws.Range("F22").Value = 0.9982129698157299
$ws.Range("I22").Value = 1.021113024687192
$ws.Range("J22").Value = 1.009310600313289
$ws.Range("K22").Value = 1.009514396634615
$ws.Range("L22").Value = 1.009555596396005
$ws.Range("M22").Value = 1.002074293599932
$ws.Range("N22").Value = 1.007305528508514

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.002722108181929
$ws.Range("D23").Value = 1.006125738468164
$ws.Range("E23").Value = 1.006154158981539
$ws.Range("F23").Value = 0.9993641489682491
$ws.Range("I23").Value = 1.021263815888672
$ws.Range("J23").Value = 1.009676388144357
$ws.Range("K23").Value = 1.009883498886546
$ws.Range("L23").Value = 1.009911804487238
$ws.Range("M23").Value = 1.003149455621163
$ws.Range("N23").Value = 1.007425702552166

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.004688196017542
$ws.Range("D24").Value = 1.007846607395211
$ws.Range("E24").Value = 1.007825411673208
$ws.Range("F24").Value = 1.003771443558863
$ws.Range("I24").Value = 1.021825465353333
$ws.Range("J24").Value = 1.011096027286714
$ws.Range("K24").Value = 1.011317947486734
$ws.Range("L24").Value = 1.011296830370928
$ws.Range("M24").Value = 1.007257964318045
$ws.Range("N24").Value = 1.007892085325282

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.006930326465401
$ws.Range("D25").Value = 1.00980886605852
$ws.Range("E25").Value = 1.009732827631846
$ws.Range("F25").Value = 1.008637385319004
$ws.Range("I25").Value = 1.022412558836226
$ws.Range("J25").Value = 1.012702720556014
$ws.Range("K25").Value = 1.012945382798696
$ws.Range("L25").Value = 1.012869597039
$ws.Range("M25").Value = 1.011777799395359
$ws.Range("N25").Value = 1.008419880898031
